$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data contained in row 2 and row 3 (the header in row 1
# and all other rows are untouched). Additionally the "substrat" info
# (AN/AO) that belonged to the Calicium adspersum record moves together
# with that record from row 2 to row 3.
#
# NOTE: this runtime's Range/Cells ".Value" GETTER is broken (it returns a
# reflection placeholder string instead of the real value), so we read and
# write through ".Formula" instead, which works correctly for numbers,
# text and booleans alike.

$row2 = 2
$row3 = 3

# Columns (by letter) whose contents differ between row 2 and row 3 and
# therefore need to be swapped. (Y/AA hold date-like text and are handled
# separately below so they are not auto-converted to date serials.)
$cols = @("A","B","D","E","F","G","H","P","Q","R","S","AW","AX","AY")

foreach ($col in $cols) {
    $addr2 = "$col$row2"
    $addr3 = "$col$row3"

    $v2 = $ws.Range($addr2).Formula
    $v3 = $ws.Range($addr3).Formula

    $ws.Range($addr2).Formula = $v3
    $ws.Range($addr3).Formula = $v2
}

# Y/AA contain plain text that looks like a date ("2005-08-16"), which
# Excel would otherwise silently reinterpret as a date serial number when
# assigned via .Formula. Force the cell to text format first, then
# restore the default "Normal" style so no stray number format sticks
# around on the cell (the source file doesn't use any custom styles).
$dateCols = @("Y","AA")
foreach ($col in $dateCols) {
    $addr2 = "$col$row2"
    $addr3 = "$col$row3"

    $v2 = $ws.Range($addr2).Formula
    $v3 = $ws.Range($addr3).Formula

    $ws.Range($addr2).NumberFormat = "@"
    $ws.Range($addr3).NumberFormat = "@"

    $ws.Range($addr2).Formula = $v3
    $ws.Range($addr3).Formula = $v2

    $ws.Range($addr2).Style = "Normal"
    $ws.Range($addr3).Style = "Normal"
}

# AN2/AO2 ("1" / "1 substratenheter # ek") move down to AN3/AO3, and are
# cleared out from row 2.
$an2 = $ws.Range("AN2").Formula
$ao2 = $ws.Range("AO2").Formula

$ws.Range("AN3").Formula = $an2
$ws.Range("AO3").Formula = $ao2

$ws.Range("AN2").Formula = ""
$ws.Range("AO2").Formula = ""
